# permirsim i katit nenkulm ne kati 1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Etazhiteti (floor structure) for the whole unit: P+NK -> P+1
$ws.Range("E3").Value = "P+1"

# Floor name: Nenkulm -> Kati-1 (both the merged description cell and the Sqarim column)
$ws.Range("F7").Value = "Kati-1"
$ws.Range("I7").Value = "Kati-1"

# Update selection to match the edited state
$ws.Range("B10:F10").Select()
